$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.226.32"
$ws.Range("E2").Value = "  +2.69%  "

$ws.Range("D3").Value = "1.916.96"
$ws.Range("E3").Value = "  +2.34%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("E4").Value = "  -0.98%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "317.38"
$ws.Range("E5").Value = "  +1.32%  "

$ws.Range("E6").Value = "  -0.99%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4840"
$ws.Range("E7").Value = "  +0.81%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3830"
$ws.Range("E8").Value = "  +1.62%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07380"
$ws.Range("E9").Value = "  -0.05%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9412"
$ws.Range("E10").Value = "  +0.07%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.88"
$ws.Range("E11").Value = "  +0.93%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07812"
$ws.Range("E12").Value = "  -0.89%  "

$ws.Range("D13").Value = "1.928.43"
$ws.Range("E13").Value = "  +2.36%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.508"
$ws.Range("E14").Value = "  +1.05%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.643"
$ws.Range("E15").Value = "  +0.54%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "91.33"
$ws.Range("E16").Value = "  +0.33%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.006"
$ws.Range("E17").Value = "  -1.04%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008844"
$ws.Range("E18").Value = "  -1.39%  "

$ws.Range("E19").Value = "  -0.86%  "

$ws.Range("D20").Value = "28.224.15"
$ws.Range("E20").Value = "  +2.54%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.87"
$ws.Range("E21").Value = "  -0.51%  "

$ws.Range("E22").Value = "  +0.35%  "

$ws.Range("D23").Value = "2.146.61"
$ws.Range("E23").Value = "  +1.82%  "

$ws.Range("E24").Value = "  +2.11%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "156.28"
$ws.Range("E25").Value = "  +1.50%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.924"
$ws.Range("E26").Value = "  -1.94%  "

$ws.Range("E27").Value = "  +0.00%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.102"
$ws.Range("E28").Value = "  +4.06%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "116.51"
$ws.Range("E29").Value = "  +0.08%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.974"
$ws.Range("E30").Value = "  -0.68%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08906"
$ws.Range("E31").Value = "  -0.30%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.359"
$ws.Range("E32").Value = "  +0.93%  "

$ws.Range("E33").Value = "  +2.66%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7760"
$ws.Range("E34").Value = "  +3.48%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.703"
$ws.Range("E35").Value = "  +2.18%  "

$ws.Range("E36").Value = "  -0.77%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02055"
$ws.Range("E37").Value = "  -0.80%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.103"
$ws.Range("E38").Value = "  -1.82%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.5540"
$ws.Range("E39").Value = "  +3.06%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05321"
$ws.Range("E40").Value = "  +0.23%  "

$ws.Range("E41").Value = "  +0.13%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.054"
$ws.Range("E42").Value = "  -0.54%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1532"
$ws.Range("E43").Value = "  +0.66%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.477"
$ws.Range("E44").Value = "  +0.44%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "10.77"
$ws.Range("E45").Value = "  +1.51%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4867"
$ws.Range("E46").Value = "  +0.43%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "106.83"
$ws.Range("E47").Value = "  +3.53%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.005"
$ws.Range("E48").Value = "  -0.99%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.660"
$ws.Range("E49").Value = "  -0.39%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "68.63"
$ws.Range("E50").Value = "  +2.15%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06116"
$ws.Range("E51").Value = "  +0.07%  "

